$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit reverts the "Expense and budget frontend" change: the sheet had
# grown from a 1-header + 1-data-row table ("Food", 1500, 2024-03-29) into a
# larger table with extra expense rows. Reverting means deleting all the rows
# that were added in between, leaving just the header row and the single
# original "Food" data row (which shifts back up into row 2).
$ws.Rows("2:7").Delete()
